$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This sheet is a flat weekly price log. The update:
#   1) inserts two brand-new rows right before the current row 629
#      (shifting the existing 629-653 block down to 631-655), and
#   2) appends one more row at the very end (new row 655) that
#      duplicates the data of what is now row 654.
# ------------------------------------------------------------------

function Set-DataRow {
    param($ws, [int]$r, [hashtable]$vals)

    $ws.Range("A$r").Value = $vals.A
    $ws.Range("B$r").Value = $vals.B
    $ws.Range("C$r").Value = $vals.C
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}

# 1) Insert two new rows above the old row 629, shifting everything else down.
$ws.Rows.Item(629).Insert()
$ws.Rows.Item(629).Insert()

# 2) Fill in the data for the two newly inserted rows (629 and 630).
$row629 = @{
    A = 4; B = "Feria Lagunitas de Puerto Montt"; C = "Los Lagos"; D = 44939
    E = 10; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102005; J = "Naranja"
    K = "Valencia"; L = "Primera"; M = 600; N = 18000; O = 19000; P = 18500
    Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 1233; T = 15
}
Set-DataRow $ws 629 $row629

$row630 = @{
    A = 4; B = "Feria Lagunitas de Puerto Montt"; C = "Los Lagos"; D = 44939
    E = 10; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102005; J = "Naranja"
    K = "Valencia"; L = "Segunda"; M = 300; N = 16000; O = 16000; P = 16000
    Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 1067; T = 15
}
Set-DataRow $ws 630 $row630

# 3) Append one more row (655) at the very end, duplicating row 654's data
#    (which itself is the shifted-down copy of the former last row, 653).
$row655 = @{
    A = 4; B = "Feria Lagunitas de Puerto Montt"; C = "Los Lagos"; D = 44335
    E = 10; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102005; J = "Naranja"
    K = "Fukumoto"; L = "Segunda"; M = 100; N = 15000; O = 15000; P = 15000
    Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 1000; T = 15
}
Set-DataRow $ws 655 $row655

# Give the new date cells (D629/D630) the same date style as the rest of column D.
$ws.Range("D629").Style = $ws.Range("D628").Style
$ws.Range("D630").Style = $ws.Range("D628").Style
$ws.Range("D655").Style = $ws.Range("D654").Style

Write-Output "done"
